# Daily attendance processing - 2025-11-15 23:20:57
# Normalize the "Recorded By" (column G) values: move a leading "System, "
# token to the end of the comma-separated list, e.g.
#   "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null) {
        $text = $val.ToString()
        if ($text.StartsWith("System, ")) {
            $rest = $text.Substring(8)
            $cell.Value = $rest + ", System"
        }
    }
}
